# Updated cryptos list on Sun Jun 23 18:22:20 UTC 2024 with GitHub Actions
# Refresh coin names/links/prices/volume deltas to the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '64.075.80'
$ws.Range("E2").Value = '  -0.30%  '

# Row 3
$ws.Range("D3").Value = '3.474.09'
$ws.Range("E3").Value = '  -0.68%  '

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").Value = "'584.99"
$ws.Range("E5").Value = '  -0.23%  '

# Row 6
$ws.Range("D6").Value = "'131.55"
$ws.Range("E6").Value = '  -1.91%  '

# Row 7
$ws.Range("E7").Value = '  +0.02%  '

# Row 8
$ws.Range("E8").Value = '  -0.95%  '

# Row 9
$ws.Range("D9").Value = "'7.66"
$ws.Range("E9").Value = '  +5.59%  '

# Row 10
$ws.Range("E10").Value = '  -1.13%  '

# Row 11
$ws.Range("E11").Value = '  -0.04%  '

# Row 12
$ws.Range("D12").Value = '4.066.84'
$ws.Range("E12").Value = '  -0.71%  '

# Row 13
$ws.Range("D13").Value = "'0.120"
$ws.Range("E13").Value = '  +0.03%  '

# Row 14
$ws.Range("D14").Value = "'0.0000178"
$ws.Range("E14").Value = '  -2.01%  '

# Row 15
$ws.Range("D15").Value = '3.476.41'
$ws.Range("E15").Value = '  -0.68%  '

# Row 16
$ws.Range("D16").Value = '64.061.20'
$ws.Range("E16").Value = '  -0.34%  '

# Row 17
$ws.Range("D17").Value = "'24.85"
$ws.Range("E17").Value = '  -3.90%  '

# Row 18
$ws.Range("D18").Value = "'10.00"
$ws.Range("E18").Value = '  +0.86%  '

# Row 19
$ws.Range("D19").Value = "'5.69"
$ws.Range("E19").Value = '  -1.03%  '

# Row 20
$ws.Range("D20").Value = "'13.41"
$ws.Range("E20").Value = '  -1.40%  '

# Row 21
$ws.Range("D21").Value = "'384.72"
$ws.Range("E21").Value = '  -2.46%  '

# Row 22
$ws.Range("D22").Value = "'0.568"
$ws.Range("E22").Value = '  -0.43%  '

# Row 23
$ws.Range("D23").Value = '3.614.14'
$ws.Range("E23").Value = '  -0.71%  '

# Row 24
$ws.Range("D24").Value = "'74.73"
$ws.Range("E24").Value = '  +0.48%  '

# Row 25
$ws.Range("E25").Value = '  +0.20%  '

# Row 26
$ws.Range("B26").Value = 'PEPE'
$ws.Range("C26").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D26").Value = "'0.0000111"
$ws.Range("E26").Value = '  -3.10%  '

# Row 27
$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = '  -0.13%  '

# Row 28
$ws.Range("B28").Value = 'PancakeSwap'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D28").Value = "'2.22"
$ws.Range("E28").Value = '  -0.57%  '

# Row 29
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Value = "'7.12"
$ws.Range("E29").Value = '  -3.67%  '

# Row 30
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").Value = "'7.96"
$ws.Range("E30").Value = '  -3.90%  '

# Row 31
$ws.Range("E31").Value = '  -4.45%  '

# Row 32
$ws.Range("B32").Value = 'Kaspa'
$ws.Range("C32").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D32").Value = "'0.153"
$ws.Range("E32").Value = '  +1.44%  '

# Row 33
$ws.Range("B33").Value = 'RenzoRestakedETH'
$ws.Range("C33").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D33").Value = '3.502.45'
$ws.Range("E33").Value = '  -0.50%  '

# Row 34
$ws.Range("B34").Value = 'USDe'
$ws.Range("C34").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = '  -0.06%  '

# Row 35
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").Value = "'22.94"
$ws.Range("E35").Value = '  -2.10%  '

# Row 36
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").Value = "'5.23"
$ws.Range("E36").Value = '  +1.43%  '

# Row 37
$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").Value = "'6.78"
$ws.Range("E37").Value = '  -1.63%  '

# Row 38
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").Value = "'1.50"
$ws.Range("E38").Value = '  -3.11%  '

# Row 39
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").Value = "'161.93"
$ws.Range("E39").Value = '  -2.19%  '

# Row 40
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").Value = "'0.0778"
$ws.Range("E40").Value = '  -0.33%  '

# Row 41
$ws.Range("B41").Value = 'Mantle'
$ws.Range("C41").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D41").Value = "'0.797"
$ws.Range("E41").Value = '  -1.05%  '

# Row 42
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = '  -0.04%  '

# Row 43
$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").Value = "'41.15"
$ws.Range("E43").Value = '  -1.87%  '

# Row 44
$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").Value = "'4.30"
$ws.Range("E44").Value = '  -2.15%  '

# Row 45
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").Value = "'1.62"
$ws.Range("E45").Value = '  -1.85%  '

# Row 46
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = "'23.73"
$ws.Range("E46").Value = '  -5.92%  '

# Row 47
$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").Value = "'1.13"
$ws.Range("E47").Value = '  -3.80%  '

# Row 48
$ws.Range("B48").Value = 'Cosmos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D48").Value = "'6.71"
$ws.Range("E48").Value = '  -0.76%  '

# Row 49
$ws.Range("B49").Value = 'SuiNetwork'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D49").Value = "'0.903"
$ws.Range("E49").Value = '  +1.15%  '

# Row 50
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '2.327.57'
$ws.Range("E50").Value = '  -5.33%  '

# Row 51
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").Value = "'0.0255"
$ws.Range("E51").Value = '  -2.57%  '
